# Commit: "Tue 24 May 2022 07:49:42 AM MSK"
#
# The document already ends with the "24 May 2022" date / horizontal-rule
# separator (an empty paragraph, the "24 May 2022 " + <hr> paragraph, and a
# trailing empty paragraph). This change appends that day's poem directly
# after that separator, i.e. at the very end of the document body:
#
#   За слово можно поплатиться
#   За слово можно и не жить
#   И мысль не будет певчей птицей
#   Вокруг других голов кружить
#   (blank line)
#   И вмажет тумаки и плюхи
#   Ладонь судьбы по головам
#   Тьмы вертухаев и рашистов
#   Пока не видно по делам
#   (blank line)
#   Безумное, шальное время
#   Сияет и не меркнет свет
#   Свободы воздух пыльной бурей
#   Еще не сперт пока у всех

$d = $word.ActiveDocument

$lines = @(
    'За слово можно поплатиться',
    'За слово можно и не жить',
    'И мысль не будет певчей птицей',
    'Вокруг других голов кружить',
    '',
    'И вмажет тумаки и плюхи',
    'Ладонь судьбы по головам',
    'Тьмы вертухаев и рашистов',
    'Пока не видно по делам',
    '',
    'Безумное, шальное время',
    'Сияет и не меркнет свет',
    'Свободы воздух пыльной бурей',
    'Еще не сперт пока у всех '
)

foreach ($line in $lines) {
    # Append a brand new paragraph at the very end of the document.
    $lastPara = $d.Paragraphs.Last
    $lastPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last

    if ($line -ne "") {
        # Set the plain text normally first.
        $newPara.Range.Text = $line

        # Re-express the paragraph's run via InsertXML, scoped to just the
        # paragraph's content (i.e. excluding the trailing paragraph mark),
        # so the xml:space="preserve" attribute is explicitly present on
        # <w:t> (matching the rest of the document) while the paragraph's
        # own <w:pPr> is left completely untouched.
        $contentStart = $newPara.Range.Start
        $contentEnd = $newPara.Range.End - 1
        $contentRng = $d.Range($contentStart, $contentEnd)

        $escapedLine = $line.Replace('&', '&amp;').Replace('<', '&lt;').Replace('>', '&gt;')
        $runXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
            '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:body><w:p><w:r><w:rPr><w:rtl w:val="0"/></w:rPr>' +
            '<w:t xml:space="preserve">' + $escapedLine + '</w:t></w:r></w:p></w:body></w:document>' +
            '</pkg:xmlData></pkg:part></pkg:package>'
        $contentRng.InsertXML($runXml)
    }
}

Write-Output "Inserted $($lines.Count) paragraphs; document now has $($d.Paragraphs.Count) paragraphs."
